# Auto-applies the numeric corrections described in the commit diff
# across the ALC, ARM, CRP, CUL, GSM, LTW, and WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 502.13
$ws.Range("I15").Value = 502.13
$ws.Range("K15").Value = 1506.39
$ws.Range("M15").Value = -1337.39

$ws.Range("H17").Value = 1672156.2
$ws.Range("J17").Value = 1672156.2
$ws.Range("L17").Value = 5016468.6
$ws.Range("N17").Value = -5016804.6

$ws.Range("H64").Value = 3842.7317
$ws.Range("I64").Value = 4066.6333
$ws.Range("J64").Value = 3232.0908
$ws.Range("K64").Value = 4066.6333
$ws.Range("L64").Value = 3232.0908
$ws.Range("M64").Value = -3818.6333
$ws.Range("N64").Value = -3728.0908

$ws.Range("H67").Value = 3842.7317
$ws.Range("I67").Value = 4066.6333
$ws.Range("J67").Value = 3232.0908
$ws.Range("K67").Value = 4066.6333
$ws.Range("L67").Value = 3232.0908
$ws.Range("M67").Value = -3208.6333
$ws.Range("N67").Value = -4948.0908

$ws.Range("H127").Value = 1248.6786
$ws.Range("I127").Value = 465.7
$ws.Range("J127").Value = 1683.6666
$ws.Range("K127").Value = 1397.1
$ws.Range("L127").Value = 5050.9998
$ws.Range("M127").Value = 3562.9
$ws.Range("N127").Value = -14970.9998

$ws.Range("H138").Value = 1684.3855
$ws.Range("I138").Value = 997.55554
$ws.Range("J138").Value = 2210.468
$ws.Range("K138").Value = 2992.66662
$ws.Range("L138").Value = 6631.404
$ws.Range("M138").Value = 2147.33338
$ws.Range("N138").Value = -16911.404

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1352889.9
$ws.Range("I122").Value = 1834556.8
$ws.Range("J122").Value = 4222.6
$ws.Range("K122").Value = 5503670.4
$ws.Range("L122").Value = 12667.8
$ws.Range("M122").Value = -5501220.4
$ws.Range("N122").Value = -17567.8

$ws.Range("H125").Value = 37799
$ws.Range("J125").Value = 37799
$ws.Range("L125").Value = 37799
$ws.Range("N125").Value = -47639

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2783.6072
$ws.Range("I31").Value = 1531.8605
$ws.Range("K31").Value = 1531.8605
$ws.Range("M31").Value = -1236.8605

$ws.Range("H34").Value = 2783.6072
$ws.Range("I34").Value = 1531.8605
$ws.Range("K34").Value = 1531.8605
$ws.Range("M34").Value = -1329.8605

$ws.Range("H86").Value = 125002470
$ws.Range("I86").Value = 250002480
$ws.Range("K86").Value = 250002480
$ws.Range("M86").Value = -250001357

$ws.Range("H89").Value = 125002470
$ws.Range("I89").Value = 250002480
$ws.Range("K89").Value = 1250012400
$ws.Range("M89").Value = -1250006784

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 223168.6
$ws.Range("I5").Value = 484.72726
$ws.Range("J5").Value = 295213.38
$ws.Range("K5").Value = 1454.18178
$ws.Range("L5").Value = 885640.14
$ws.Range("M5").Value = -1342.18178
$ws.Range("N5").Value = -885864.14

$ws.Range("H63").Value = 3674.3333
$ws.Range("I63").Value = 2011.5
$ws.Range("J63").Value = 7000
$ws.Range("K63").Value = 6034.5
$ws.Range("L63").Value = 21000
$ws.Range("M63").Value = -5285.5
$ws.Range("N63").Value = -22498

$ws.Range("H64").Value = 3484.8
$ws.Range("J64").Value = 4133.3335
$ws.Range("L64").Value = 12400.0005
$ws.Range("N64").Value = -12940.0005

$ws.Range("H66").Value = 3674.3333
$ws.Range("I66").Value = 2011.5
$ws.Range("J66").Value = 7000
$ws.Range("K66").Value = 18103.5
$ws.Range("L66").Value = 63000
$ws.Range("M66").Value = -14359.5
$ws.Range("N66").Value = -70488

$ws.Range("H67").Value = 3484.8
$ws.Range("J67").Value = 4133.3335
$ws.Range("L67").Value = 12400.0005
$ws.Range("N67").Value = -14272.0005

$ws.Range("H69").Value = 1843.2727
$ws.Range("I69").Value = 930.6667
$ws.Range("J69").Value = 2185.5
$ws.Range("K69").Value = 2792.0001
$ws.Range("L69").Value = 6556.5
$ws.Range("M69").Value = -1981.0001
$ws.Range("N69").Value = -8178.5

$ws.Range("H70").Value = 2702.6428
$ws.Range("J70").Value = 3823.111
$ws.Range("L70").Value = 11469.333
$ws.Range("N70").Value = -12099.333

$ws.Range("H72").Value = 1843.2727
$ws.Range("I72").Value = 930.6667
$ws.Range("J72").Value = 2185.5
$ws.Range("K72").Value = 8376.0003
$ws.Range("L72").Value = 19669.5
$ws.Range("M72").Value = -4320.0003
$ws.Range("N72").Value = -27781.5

$ws.Range("H73").Value = 2702.6428
$ws.Range("J73").Value = 3823.111
$ws.Range("L73").Value = 11469.333
$ws.Range("N73").Value = -13653.333

$ws.Range("H76").Value = 2738.3333
$ws.Range("I76").Value = 1000
$ws.Range("J76").Value = 3607.5
$ws.Range("K76").Value = 3000
$ws.Range("L76").Value = 10822.5
$ws.Range("M76").Value = -2617
$ws.Range("N76").Value = -11588.5

$ws.Range("H79").Value = 2738.3333
$ws.Range("I79").Value = 1000
$ws.Range("J79").Value = 3607.5
$ws.Range("K79").Value = 3000
$ws.Range("L79").Value = 10822.5
$ws.Range("M79").Value = -1674
$ws.Range("N79").Value = -13474.5

$ws.Range("H88").Value = 3091.6
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 3364.5
$ws.Range("K88").Value = 6000
$ws.Range("L88").Value = 10093.5
$ws.Range("M88").Value = -5572
$ws.Range("N88").Value = -10949.5

$ws.Range("H91").Value = 3091.6
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 3364.5
$ws.Range("K91").Value = 6000
$ws.Range("L91").Value = 10093.5
$ws.Range("M91").Value = -4518
$ws.Range("N91").Value = -13057.5

$ws.Range("H94").Value = 2810.3333
$ws.Range("I94").Value = 1262
$ws.Range("J94").Value = 3120
$ws.Range("K94").Value = 3786
$ws.Range("L94").Value = 9360
$ws.Range("M94").Value = -3110
$ws.Range("N94").Value = -10712

$ws.Range("H100").Value = 3720
$ws.Range("I100").Value = 3780
$ws.Range("J100").Value = 3700
$ws.Range("K100").Value = 11340
$ws.Range("L100").Value = 11100
$ws.Range("M100").Value = -10529
$ws.Range("N100").Value = -12722

$ws.Range("H103").Value = 5486.091
$ws.Range("I103").Value = 7078.143
$ws.Range("J103").Value = 2700
$ws.Range("K103").Value = 21234.429
$ws.Range("L103").Value = 8100
$ws.Range("M103").Value = -20355.429
$ws.Range("N103").Value = -9858

$ws.Range("H122").Value = 8064.2856
$ws.Range("I122").Value = 440.6
$ws.Range("J122").Value = 27123.5
$ws.Range("K122").Value = 3965.4
$ws.Range("L122").Value = 244111.5
$ws.Range("M122").Value = -1515.4
$ws.Range("N122").Value = -249011.5

$ws.Range("H135").Value = 223168.6
$ws.Range("I135").Value = 484.72726
$ws.Range("J135").Value = 295213.38
$ws.Range("K135").Value = 4362.54534
$ws.Range("L135").Value = 2656920.42
$ws.Range("M135").Value = -1827.54534
$ws.Range("N135").Value = -2661990.42

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5222.6284
$ws.Range("I70").Value = 5345.865
$ws.Range("J70").Value = 4866.6113
$ws.Range("K70").Value = 5345.865
$ws.Range("L70").Value = 4866.6113
$ws.Range("M70").Value = -5075.865
$ws.Range("N70").Value = -5406.6113

$ws.Range("H73").Value = 5222.6284
$ws.Range("I73").Value = 5345.865
$ws.Range("J73").Value = 4866.6113
$ws.Range("K73").Value = 5345.865
$ws.Range("L73").Value = 4866.6113
$ws.Range("M73").Value = -4409.865
$ws.Range("N73").Value = -6738.6113

$ws.Range("H102").Value = 893729.5
$ws.Range("I102").Value = 1695765.6
$ws.Range("J102").Value = 2578.2222
$ws.Range("K102").Value = 1695765.6
$ws.Range("L102").Value = 2578.2222
$ws.Range("M102").Value = -1694143.6
$ws.Range("N102").Value = -5822.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 920.7895
$ws.Range("I46").Value = 719.6
$ws.Range("J46").Value = 992.6429000000001
$ws.Range("K46").Value = 719.6
$ws.Range("L46").Value = 992.6429000000001
$ws.Range("M46").Value = -531.6
$ws.Range("N46").Value = -1368.6429

$ws.Range("H82").Value = 22688500
$ws.Range("I82").Value = 5001000
$ws.Range("K82").Value = 5001000
$ws.Range("M82").Value = -5000639

$ws.Range("H85").Value = 22688500
$ws.Range("I85").Value = 5001000
$ws.Range("K85").Value = 5001000
$ws.Range("M85").Value = -4999752

$ws.Range("H100").Value = 1464.375
$ws.Range("I100").Value = 1387.8572
$ws.Range("K100").Value = 1387.8572
$ws.Range("M100").Value = -846.8571999999999

$ws.Range("H122").Value = 3881427
$ws.Range("I122").Value = 4765698
$ws.Range("K122").Value = 14297094
$ws.Range("M122").Value = -14294644

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2187.0454
$ws.Range("I132").Value = 1394.5
$ws.Range("J132").Value = 4300.5
$ws.Range("K132").Value = 4183.5
$ws.Range("L132").Value = 12901.5
$ws.Range("M132").Value = -1653.5
$ws.Range("N132").Value = -17961.5

"done"
